$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 'maa://24702 (94.27), maa://25390 (95.93), maa://36681 (87.01)'
$ws.Range("L2").Value = '*maa://24633 (56.05), *maa://30515 (69.31), *maa://34787 (72.86), ***maa://20792 (11.93), maa://39402 (88.37), ***maa://29083 (27.78)'
$ws.Range("AF2").Value = 'maa://25251 (91.75), ***maa://21730 (22.86), ***maa://39501 (18.18), *maa://36675 (60.0)'
$ws.Range("H3").Value = 'maa://21247 (98.45), *maa://22748 (60.0)'
$ws.Range("P3").Value = 'maa://21249 (94.22), maa://26254 (96.0)'
$ws.Range("D4").Value = 'maa://24632 (93.38), **maa://24303 (33.33), maa://22499 (86.67), maa://22746 (100.0)'
$ws.Range("L7").Value = 'maa://28624 (93.18), maa://24957 (97.62)'
$ws.Range("AF7").Value = '*maa://26191 (68.35), *maa://36671 (69.39), *maa://42530 (64.29)'
$ws.Range("A8").Value = '更新日期：2024.12.31 13:17:57'
$ws.Range("X8").Value = 'maa://21411 (95.93)'
$ws.Range("D10").Value = '***maa://25695 (19.44), **maa://32237 (40.48), ***maa://34206 (21.74), ***maa://39951 (16.67), ***maa://39243 (28.57)'
$ws.Range("H12").Value = 'maa://21867 (89.76)'
$ws.Range("AB12").Value = 'maa://23669 (95.41), maa://36677 (92.31), maa://39872 (90.0)'
$ws.Range("P13").Value = 'maa://22676 (91.96), *maa://22583 (75.0), *maa://22500 (57.78)'
$ws.Range("X13").Value = '*maa://34957 (79.69), *maa://22768 (51.61)'
$ws.Range("T15").Value = 'maa://23892 (97.4)'
$ws.Range("T16").Value = 'maa://22729 (95.39), *maa://28648 (68.85), maa://36674 (82.93)'
$ws.Range("H18").Value = 'maa://24421 (90.5)'
$ws.Range("L18").Value = 'maa://22466 (88.97), *maa://22732 (50.6)'
$ws.Range("D20").Value = 'maa://21432 (89.86), maa://25198 (93.14), *maa://20795 (51.18), maa://36680 (93.55)'
$ws.Range("L20").Value = 'maa://41331 (85.15)'
$ws.Range("AB21").Value = 'maa://21443 (80.17), ***maa://23820 (29.82)'
$ws.Range("L22").Value = 'maa://27127 (85.44), *maa://22751 (73.85)'
$ws.Range("X22").Value = 'maa://21282 (98.47), *maa://37649 (66.67)'
$ws.Range("L23").Value = 'maa://39756 (94.3), maa://39875 (93.75)'
$ws.Range("X24").Value = 'maa://29988 (86.84), maa://23504 (93.15), **maa://22892 (39.58), *maa://25141 (76.98), *maa://36663 (78.26), ***maa://22815 (23.08)'
$ws.Range("AF24").Value = 'maa://22523 (85.71), maa://36672 (80.77), maa://29910 (92.59), **maa://21440 (34.55)'
$ws.Range("X25").Value = '*maa://29890 (76.74)'
$ws.Range("AF25").Value = 'maa://20108 (96.27), maa://24621 (96.58), maa://36676 (96.77), maa://22771 (85.71), maa://37772 (100.0)'
$ws.Range("X28").Value = 'maa://39929 (89.97), ***maa://39723 (14.29), maa://41749 (91.38)'
$ws.Range("AF28").Value = 'maa://36660 (92.99), *maa://36701 (64.29)'
$ws.Range("L29").Value = 'maa://28432 (92.9), *maa://28440 (76.84), maa://31400 (100.0), *maa://28650 (71.43)'
$ws.Range("AB30").Value = 'maa://42979 (96.4), maa://45045 (100.0)'
$ws.Range("L31").Value = 'maa://35926 (93.66), maa://36258 (83.67), *maa://43904 (77.78)'
$ws.Range("H32").Value = 'maa://21895 (97.3), maa://36667 (98.39), **maa://20793 (38.78), maa://22760 (100.0)'
$ws.Range("T32").Value = 'maa://42859 (96.34), maa://41108 (87.76), maa://41238 (96.3)'
$ws.Range("L35").Value = 'maa://41296 (96.69)'
$ws.Range("T35").Value = 'maa://24842 (94.12)'
$ws.Range("H39").Value = 'maa://25199 (84.82), maa://36670 (87.64), maa://30434 (89.39), ***maa://25036 (16.0), *maa://44165 (66.67), maa://45059 (100.0)'
$ws.Range("H44").Value = 'maa://29768 (97.84), maa://27728 (96.0)'
$ws.Range("H46").Value = 'maa://35931 (92.61), maa://43901 (88.89)'
$ws.Range("H47").Value = 'maa://27410 (96.2), maa://29661 (97.86), maa://28038 (84.62)'
$ws.Range("H53").Value = 'maa://32534 (93.56), **maa://32434 (34.78)'
$ws.Range("H55").Value = 'maa://32532 (91.92)'
$ws.Range("H57").Value = 'maa://25176 (98.25)'
